# Refactor + Excel improvement + Usability improvements
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Create")

# Update row 2: B2 becomes text "a" instead of numeric 1, D2 gets cleared (was duplicate "sp1")
$ws.Cells.Item(2, 2).Value = "a"
$ws.Range("D2").ClearContents()

# Remove row 3 (fs2 / 2 / sp2) entirely -- clear contents instead of deleting the
# row so the data validation ranges (which went down to row 1000) are not shifted.
$ws.Range("A3:C3").ClearContents()

# Update selection to reflect new active cell
$ws.Range("C2").Select() | Out-Null
